$d = $word.ActiveDocument

# 1. Title / heading text changes (appears twice: Heading1 and bold run near end)
#    Replace:=2 (wdReplaceAll) replaces every occurrence in the document in one call.
$d.Content.Find.Execute("Play Leprechaun Goes Wild Free| Exciting Irish-themed Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play Leprechaun Goes Wild for Free", 2)

# 2. "What we like" bullet list changes
$d.Content.Find.Execute("Fun and exciting bonus feature", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting bonus feature with free spins", 2)

$d.Content.Find.Execute("Detailed and colorful graphics", $true, $false, $false, $false, $false, $true, 1, $false, "Visually appealing graphics and details", 2)

$d.Content.Find.Execute("Lively and engaging sound effects", $true, $false, $false, $false, $false, $true, 1, $false, "Engaging and upbeat sound effects", 2)

# 3. "What we don't like" bullet list changes
$d.Content.Find.Execute("Limited variety in symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Limited variety of symbols", 2)

$d.Content.Find.Execute("May be too similar to other Irish-themed slots", $true, $false, $false, $false, $false, $true, 1, $false, "Theme may not appeal to everyone", 2)

# 4. Meta description (italic run)
$d.Content.Find.Execute("Read our review of Leprechaun Goes Wild slot game. Play for free and enjoy the Irish folklore theme, bonus features, and high payout potential.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Leprechaun Goes Wild slot game and play for free.", 2)
